$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 282, shifting existing rows 282:301 down to 283:302
$ws.Rows.Item(282).Insert()

# Fill in the new row 282 with data (mirrors the surrounding Repollo rows)
$ws.Cells.Item(282, 1).Value = 5
$ws.Cells.Item(282, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(282, 3).Value = "Maule"
$ws.Cells.Item(282, 4).Value = 44746
$ws.Cells.Item(282, 5).Value = 7
$ws.Cells.Item(282, 6).Value = 100112006
$ws.Cells.Item(282, 7).Value = "Repollo"
$ws.Cells.Item(282, 8).Value = "Crespo record"
$ws.Cells.Item(282, 9).Value = "Primera"
$ws.Cells.Item(282, 10).Value = 3000
$ws.Cells.Item(282, 11).Value = 1100
$ws.Cells.Item(282, 12).Value = 1100
$ws.Cells.Item(282, 13).Value = 1100
$ws.Cells.Item(282, 14).Value = "`$/unidad"
$ws.Cells.Item(282, 15).Value = "Región del Maule"
$ws.Cells.Item(282, 16).Value = 1100
$ws.Cells.Item(282, 17).Value = 1
$ws.Cells.Item(282, 18).Value = "Hortaliza"
